# location bulk issue resolved
# Replace the "skills" list (Python/JavaScript/React Js/css/html/angular)
# with a short "location" list (delhi/jaipur J) and drop the now-unused rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove rows 4-7 entirely (shifts nothing below them up; they were the last rows)
$ws.Range("A4:A7").EntireRow.Delete()

# Update remaining data rows with the new location values
$ws.Range("A2").Value = "delhi"
$ws.Range("A3").Value = "jaipur J"

# Row 3 loses its special styling/height (it becomes the new "last row",
# matching the previously-unstyled last row's look-and-feel)
$ws.Range("A3").Style = "Normal"
$ws.Rows(3).AutoFit()

# Move the active selection up to A4 (the new first empty row after data)
$ws.Range("A4").Select()
